$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Modify existing cells ---
$ws.Cells.Item(3, 13).Value = 40
$ws.Cells.Item(22, 2).Value = "An experiment at the non-interacting lattice stage."
$ws.Cells.Item(22, 7).Value = "DensityFit;AtomNumber;CenterFit"
$ws.Cells.Item(22, 13).Value = 0.20000000000000001
$ws.Cells.Item(22, 14).Value = "LinearFit1D"
$ws.Cells.Item(23, 2).Value = "An experiment at the non-interacting BEC stage."
$ws.Cells.Item(26, 2).Value = "An experiment at the non-interacting lattice stage. "
$ws.Cells.Item(26, 7).Value = "DensityFit;CenterFit;AtomNumber"

# --- Add new rows ---
# Row 27: partialEvapA
$ws.Cells.Item(27, 1).Value = "partialEvapA"
$ws.Cells.Item(27, 2).Value = "An experiment at evaporation stage D."
$ws.Cells.Item(27, 3).Value = "TOP"
$ws.Cells.Item(27, 4).Value = "EvapDOdt1"
$ws.Cells.Item(27, 5).Value = 4
$ws.Cells.Item(27, 6).Value = "RunIndex"
$ws.Cells.Item(27, 7).Value = "DensityFit;AtomNumber;CenterFit"
$ws.Cells.Item(27, 8).Value = "LSR"
$ws.Cells.Item(27, 9).Value = "LF"
$ws.Cells.Item(27, 10).Value = "RandomPolarization"
$ws.Cells.Item(27, 11).Value = 8
$ws.Cells.Item(27, 12).Value = "BosonicGaussianFit1D"
$ws.Cells.Item(27, 13).Value = 30
$ws.Cells.Item(27, 14).Value = "LinearFit1D"

# Row 28: partialEvapB
$ws.Cells.Item(28, 1).Value = "partialEvapB"
$ws.Cells.Item(28, 2).Value = "An experiment at evaporation stage A."
$ws.Cells.Item(28, 3).Value = "TOP"
$ws.Cells.Item(28, 4).Value = "EvapDOdt1"
$ws.Cells.Item(28, 5).Value = 4
$ws.Cells.Item(28, 6).Value = "RunIndex"
$ws.Cells.Item(28, 7).Value = "DensityFit;AtomNumber;CenterFit"
$ws.Cells.Item(28, 8).Value = "LSR"
$ws.Cells.Item(28, 9).Value = "LF"
$ws.Cells.Item(28, 10).Value = "RandomPolarization"
$ws.Cells.Item(28, 11).Value = 8
$ws.Cells.Item(28, 12).Value = "BosonicGaussianFit1D"
$ws.Cells.Item(28, 13).Value = 30
$ws.Cells.Item(28, 14).Value = "LinearFit1D"

# Row 29: partialEvapC
$ws.Cells.Item(29, 1).Value = "partialEvapC"
$ws.Cells.Item(29, 2).Value = "An experiment at evaporation stage B."
$ws.Cells.Item(29, 3).Value = "TOP"
$ws.Cells.Item(29, 4).Value = "EvapDOdt1"
$ws.Cells.Item(29, 5).Value = 4
$ws.Cells.Item(29, 6).Value = "RunIndex"
$ws.Cells.Item(29, 7).Value = "DensityFit;AtomNumber;CenterFit"
$ws.Cells.Item(29, 8).Value = "LSR"
$ws.Cells.Item(29, 9).Value = "LF"
$ws.Cells.Item(29, 10).Value = "RandomPolarization"
$ws.Cells.Item(29, 11).Value = 8
$ws.Cells.Item(29, 12).Value = "BosonicGaussianFit1D"
$ws.Cells.Item(29, 13).Value = 30
$ws.Cells.Item(29, 14).Value = "LinearFit1D"

# Row 30: NiBecCameraSBB
$ws.Cells.Item(30, 1).Value = "NiBecCameraSBB"
$ws.Cells.Item(30, 2).Value = "An experiment at the non-interacting BEC stage."
$ws.Cells.Item(30, 3).Value = "SBB"
$ws.Cells.Item(30, 4).Value = "Full"
$ws.Cells.Item(30, 5).Value = 4
$ws.Cells.Item(30, 6).Value = "RunIndex"
$ws.Cells.Item(30, 7).Value = "DensityFit;AtomNumber;CenterFit"
$ws.Cells.Item(30, 8).Value = "LSR"
$ws.Cells.Item(30, 9).Value = "NI"
$ws.Cells.Item(30, 10).Value = "StrongLight"
$ws.Cells.Item(30, 11).Value = 8
$ws.Cells.Item(30, 12).Value = "BosonicGaussianFit1D"
$ws.Cells.Item(30, 13).Value = 0.20000000000000001
$ws.Cells.Item(30, 14).Value = "ParabolicFit1D"

# Row 31: NiLatticeXvNi
$ws.Cells.Item(31, 1).Value = "NiLatticeXvNi"
$ws.Cells.Item(31, 2).Value = "An experiment at the non-interacting lattice stage. Scan XV_NI"
$ws.Cells.Item(31, 3).Value = "TOP"
$ws.Cells.Item(31, 4).Value = "NiLattice"
$ws.Cells.Item(31, 5).Value = 4
$ws.Cells.Item(31, 6).Value = "XV_NI"
$ws.Cells.Item(31, 7).Value = "DensityFit;AtomNumber;CenterFit"
$ws.Cells.Item(31, 8).Value = "LSR"
$ws.Cells.Item(31, 9).Value = "HF"
$ws.Cells.Item(31, 10).Value = "StrongLight"
$ws.Cells.Item(31, 11).Value = 8
$ws.Cells.Item(31, 12).Value = "BosonicGaussianFit1D"
$ws.Cells.Item(31, 13).Value = 0.20000000000000001
$ws.Cells.Item(31, 14).Value = "LinearFit1D"

# Row 32: NiLatticeKdPulse
$ws.Cells.Item(32, 1).Value = "NiLatticeKdPulse"
$ws.Cells.Item(32, 2).Value = "An experiment at the non-interacting lattice stage. Pulsed on Kapitza Dirac."
$ws.Cells.Item(32, 3).Value = "TOP"
$ws.Cells.Item(32, 4).Value = "NiLattice"
$ws.Cells.Item(32, 5).Value = 4
$ws.Cells.Item(32, 6).Value = "RunIndex"
$ws.Cells.Item(32, 7).Value = "KapitzaDirac"
$ws.Cells.Item(32, 8).Value = "LSR"
$ws.Cells.Item(32, 9).Value = "HF"
$ws.Cells.Item(32, 10).Value = "StrongLight"
$ws.Cells.Item(32, 11).Value = 8
$ws.Cells.Item(32, 12).Value = "BosonicGaussianFit1D"
$ws.Cells.Item(32, 13).Value = 0.20000000000000001
$ws.Cells.Item(32, 14).Value = "SineFit1D"

# Row 33: NiLatticeBo
$ws.Cells.Item(33, 1).Value = "NiLatticeBo"
$ws.Cells.Item(33, 2).Value = "A Bloch oscillation experiment at the non-interacting lattice stage."
$ws.Cells.Item(33, 3).Value = "TOP"
$ws.Cells.Item(33, 4).Value = "NiLattice"
$ws.Cells.Item(33, 5).Value = 4
$ws.Cells.Item(33, 6).Value = "latticehold"
$ws.Cells.Item(33, 7).Value = "DensityFit;AtomNumber;CenterFit"
$ws.Cells.Item(33, 8).Value = "LSR"
$ws.Cells.Item(33, 9).Value = "HF"
$ws.Cells.Item(33, 10).Value = "StrongLight"
$ws.Cells.Item(33, 11).Value = 8
$ws.Cells.Item(33, 12).Value = "BosonicGaussianFit1D"
$ws.Cells.Item(33, 13).Value = 0.20000000000000001
$ws.Cells.Item(33, 14).Value = "SineFit1D"

# Row 34: NiLatticeTransfer
$ws.Cells.Item(34, 1).Value = "NiLatticeTransfer"
$ws.Cells.Item(34, 2).Value = "An experiment at the non-interacting lattice stage."
$ws.Cells.Item(34, 3).Value = "TOP"
$ws.Cells.Item(34, 4).Value = "NiLattice"
$ws.Cells.Item(34, 5).Value = 4
$ws.Cells.Item(34, 6).Value = "dummy2"
$ws.Cells.Item(34, 7).Value = "AtomNumber;DensityFit"
$ws.Cells.Item(34, 8).Value = "LSR"
$ws.Cells.Item(34, 9).Value = "HF"
$ws.Cells.Item(34, 10).Value = "StrongLight"
$ws.Cells.Item(34, 11).Value = 8
$ws.Cells.Item(34, 12).Value = "BosonicGaussianFit1D"
$ws.Cells.Item(34, 13).Value = 0.20000000000000001
$ws.Cells.Item(34, 14).Value = "LinearFit1D"
